# Daily attendance processing - 2026-01-27 06:09:16
# Normalize the ordering of names/emails in the "Recorded By" (column G) list
# for rows whose recorder list was stored in a stale/inconsistent order.
# The fix re-orders the last two entries of the comma-separated "Recorded By"
# value for the specific known-stale combinations produced by the attendance
# recorder merge.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-match lookup of stale "Recorded By" text -> corrected text.
$map = @{
    "dnasr281@gmail.com, System"            = "System, dnasr281@gmail.com"
    "backup@backdoor.com, system, System"   = "backup@backdoor.com, System, system"
    "dnasr281@gmail.com, admin@admin.com"   = "admin@admin.com, dnasr281@gmail.com"
}

$lastRow = $ws.UsedRange.Rows.Count
$recordedByCol = 7  # Column G = "Recorded By"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $recordedByCol)
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
